# Update schedule and EC
# - Row 8 (week 7) homework cell (D8) gets expanded from "HW6" to the full
#   assignment description, and a matching note (E8) is added, mirroring
#   the HW4/HW5 rows above it.
# - The active selection moves from D8 to C8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D8").Value = "HW6: tie a clove hitch and cleat hitch"
$ws.Range("E8").Value = "bring your line to class; be prepared to go outside"

$ws.Range("C8").Select()
